$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 215; this pushes the existing rows
# 215..322 down to 216..323 (dimension grows from R322 to R323).
$ws.Rows.Item(215).Insert()

# Populate the newly inserted row 215 with the new data record.
$ws.Cells.Item(215, 1).Value = 7
$ws.Cells.Item(215, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(215, 3).Value = "Ñuble"
$ws.Cells.Item(215, 4).Value = 45202
$ws.Cells.Item(215, 5).Value = 16
$ws.Cells.Item(215, 6).Value = 100112040
$ws.Cells.Item(215, 7).Value = "Cilantro"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 280
$ws.Cells.Item(215, 11).Value = 1500
$ws.Cells.Item(215, 12).Value = 1500
$ws.Cells.Item(215, 13).Value = 1500
$ws.Cells.Item(215, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(215, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(215, 16).Value = 1500
$ws.Cells.Item(215, 17).Value = 1
$ws.Cells.Item(215, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# rest of column D.
$ws.Cells.Item(215, 4).NumberFormat = $ws.Cells.Item(216, 4).NumberFormat
